# Remove the explicit transition speed ("spd") override on the slides whose
# transition was set to slow/medium, restoring the (unspecified) default speed
# while leaving the rest of the transition (advance-on-click / advance-on-time,
# and the morph/fade effect) untouched.
$p = $ppt.ActivePresentation

$slideIndexes = @(2, 3, 4, 5)

foreach ($idx in $slideIndexes) {
    $slide = $p.Slides.Item($idx)
    $slide.SlideShowTransition.Speed = $null
}
